$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") rows 2-24: update date serial value from 45227 to 45228
$ws.Range("C2:C24").Value = 45228
